# Add a data-quality note in the "dataText" column (E) for the local-skills
# rows (2-9) of the dataText sheet, per commit:
# "Add in note on the local skills page about the data quality"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dataText")

# Build the note text (HTML snippet stored as plain text, with a literal
# right single quotation mark in "ONS' view").
$rightQuote = [char]0x2019
$note = "<p>" + "`n" +
  "  ONS published a " + "`n" +
  "  <a href=`"https://osr.statisticsauthority.gov.uk/correspondence/michael-keoghan-to-siobhan-tuohy-smith-request-to-suspend-aps-accreditation/`">response to OSR</a> about the current quality of Annual Population Survey (APS) (and Labour Force Survey) outputs. ONS asked OSR to temporarily suspend accreditation of all APS-based ONS outputs. There has since been a " + "`n" +
  "  <a href=`"https://osr.statisticsauthority.gov.uk/correspondence/ed-humpherson-to-michael-keoghan-suspension-of-the-accredited-official-statistics-status-for-the-estimates-ons-produces-from-the-annual-population-survey/`">response letter from OSR</a>. Overall, ONS" + $rightQuote + " view on the quality of the APS is that while it is robust for national and headline regional estimates, there are concerns with the quality of estimates for smaller segments of the population, such as local authority geographies. ONS will publish an explanatory note later this year providing guidance to users on the quality of current APS and will be used to inform further work ONS is undertaking to improve quality of the survey." + "`n" +
  "</p>"

# Write the note into E2:E9 (the dataText column for each local-skills metric row)
$range = $ws.Range("E2:E9")
$range.Value = $note

# Match the new cell style: left/centre aligned, wrapped text
$range.WrapText = $true
$range.HorizontalAlignment = -4131   # xlLeft
$range.VerticalAlignment = -4108     # xlCenter

# Grow the rows to fit the much longer text (Excel's maximum row height)
foreach ($r in 2..9) {
    $ws.Rows.Item($r).RowHeight = 409.6
}

# Reflect the scrolled/selected view used while reviewing the new note
$ws.Activate()
$excel.Goto($ws.Range("A10"), $true)
$ws.Range("E10").Select()
